$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 153, pushing existing rows 153:164 down to 154:165
$ws.Rows.Item(153).Insert()

# Populate the new row 153 with the new record
$ws.Cells.Item(153, 1).Value = 10
$ws.Cells.Item(153, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(153, 3).Value = "La Araucanía"
$ws.Cells.Item(153, 4).Value = 44769
$ws.Cells.Item(153, 4).Style = $ws.Cells.Item(154, 4).Style
$ws.Cells.Item(153, 4).NumberFormat = $ws.Cells.Item(154, 4).NumberFormat
$ws.Cells.Item(153, 5).Value = 9
$ws.Cells.Item(153, 6).Value = 100112012
$ws.Cells.Item(153, 7).Value = "Espinaca"
$ws.Cells.Item(153, 8).Value = "Sin especificar"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 80
$ws.Cells.Item(153, 11).Value = 9000
$ws.Cells.Item(153, 12).Value = 9000
$ws.Cells.Item(153, 13).Value = 9000
$ws.Cells.Item(153, 14).Value = "$/docena de atados"
$ws.Cells.Item(153, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(153, 16).Value = 3000
$ws.Cells.Item(153, 17).Value = 3
$ws.Cells.Item(153, 18).Value = "Hortaliza"
